# Generate Report for handback
#
# For each localized-language sheet ("zh-cn" and "de-de") this:
#   - flips the Status column (B) from "Ready for handoff" to
#     "Handed back: in sync with en-US" for the two data rows,
#   - fills in "Latest Target File" (E) / "Latest Handback File" (F) with
#     the same file links already recorded in "Source File Name" (A) and
#     "Latest Handoff File" (C),
#   - stamps "Latest Handback DateTime" (G) with the handback timestamp, and
#   - flips "Handoff Reason" (H) from "Ignored" to "Include".

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

function Get-HyperlinkForRange($ws, $rng) {
    $target = $rng.Address()
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $target) {
            return $h
        }
    }
    return $null
}

function Update-HandbackSheet($SheetName, $HandbackDateTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($row in 2, 3) {

        # Status -> handed back
        $ws.Cells.Item($row, 2).Value = $statusHandedBack

        # Source md link (column A) mirrored into "Latest Target File" (E)
        $srcRange = $ws.Cells.Item($row, 1)
        $srcLink = Get-HyperlinkForRange $ws $srcRange
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $srcLink.Address, "", "", $srcRange.Value()) | Out-Null

        # Handoff xlf link (column C) mirrored into "Latest Handback File" (F)
        $handoffRange = $ws.Cells.Item($row, 3)
        $handoffLink = Get-HyperlinkForRange $ws $handoffRange
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $handoffLink.Address, "", "", $handoffRange.Value()) | Out-Null

        # Latest Handback DateTime (G)
        $ws.Cells.Item($row, 7).Value = $HandbackDateTime

        # Handoff Reason (H) -> Include
        $ws.Cells.Item($row, 8).Value = "Include"
    }
}

Update-HandbackSheet "zh-cn" "2016-01-11 13:00:35"
Update-HandbackSheet "de-de" "2016-01-11 13:01:06"

# The "Overview" sheet mirrors the per-language status for each file, so it
# also needs to reflect the new handed-back status (column B = zh-cn status,
# column C = de-de status) for the two tracked files (rows 2 and 3).
$overview = $wb.Worksheets.Item("Overview")
foreach ($row in 2, 3) {
    $overview.Cells.Item($row, 2).Value = $statusHandedBack
    $overview.Cells.Item($row, 3).Value = $statusHandedBack
}
